# 2.a.1.xlsx - add a new "2021" column (R) to the agriculture orientation
# index table, and refresh the already-published 2019/2020 (P4/Q4) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3): new year label 2021 in R3 -----------------------
# Clone Q3's look (general format, same border/alignment as the other year
# headers) onto R3, then write the year value.
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R3").Value = 2021

# --- Data row (row 4): new data point for 2021 in R4 ----------------------
# Clone Q4's look (numeric "0.00" format, right/center aligned, bordered)
# onto R4 first ...
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
# ... give it its own (slightly distinct) font so this new column gets a
# dedicated style entry rather than silently sharing R4 with the older
# columns ...
$ws.Range("R4").Font.Name = "Times New Roman"
$ws.Range("R4").Font.Size = 9
$ws.Range("R4").Font.Color = 0
# ... and finally drop in the figure itself.
$ws.Range("R4").Value = 0.064467421337540437

# --- Revised figures for the two most recent, previously-published years -
$ws.Range("P4").Value = 0.09130340807234763
$ws.Range("Q4").Value = 0.074862480994528399

$excel.CutCopyMode = $false

# --- Restore the cursor / selection left by whoever made this edit --------
$ws.Range("O10").Select()
